$p = $ppt.ActivePresentation

# --- Update slide 1 (title slide) subtitle text ---
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(2).TextFrame.TextRange.Text = "BitBlocks `r `rPolicy Framework `r `rPurpose `rThis policy framework provides a clear guide for ethical conduct across all organizational `rlevels of BitBlocks. `r `rIntroduction `r➢​ Section 1 outlines foundational principles of integrity, respect, and professionalism"

# --- Update existing slides 2-8 content/titles ---
$s = $p.Slides.Item(2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Key Points"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "The BitBlocks Policy Framework provides a clear guide for ethical conduct across all organizational levels of BitBlocks`rThe policy framework outlines foundational principles of integrity, respect, and professionalism"

$s = $p.Slides.Item(3)
$s.Shapes.Item(2).TextFrame.TextRange.Text = "The Policy Change Procedure ensures that modifications to the BitBlocks policy are review ."

$s = $p.Slides.Item(4)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Key Points"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "All employees must act with honesty, transparency, and fairness in all professional interactions`rTreat colleagues, clients,janitorial staff, security team and stakeholders with respect"

$s = $p.Slides.Item(5)
$s.Shapes.Item(1).TextFrame.TextRange.Text = ""
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Protect sensitive company, client, and user data`rAdhere to all applicable laws, regulations, and industry standards ."

$s = $p.Slides.Item(6)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Key Points"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "CEO/Board-Level Policies include vision and leadership, accountability, vision and mission`rEthical Oversight: Ethical oversight, ensure compliance with the ACM and IEEE codes"

$s = $p.Slides.Item(7)
$s.Shapes.Item(1).TextFrame.TextRange.Text = ""
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Accountability: Accountability:  Prioritize employees' well-being and mental health ."

$s = $p.Slides.Item(8)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Key Points"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Team Leadership: Foster a collaborative and inclusive work environment where team members feel valued and supported`rProject Management:  Ensure projects are delivered on time, within budget, and meet quality standards"

# --- Add new slides 9-24 (Title and Content layout) ---
$s = $p.Slides.Add(9, 2)
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Stakeholder Engagement:  Engage with stakeholders (employees, clients, investors, and the public)"

$s = $p.Slides.Add(10, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Key Points"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Code of Conduct: Adhere to the company’s code of conduct and ethical guidelines in all professional activities`rCollaboration:  Collaboration with team members, sharing knowledge and sharing knowledge"

$s = $p.Slides.Add(11, 2)
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Professional Development:  Support the professional growth of team members through training,  career development opportunities ."

$s = $p.Slides.Add(12, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Key Points"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Employees are entitled to 20 paid annual leave days per year`rSalaries are disbursed on the 5th of every month"

$s = $p.Slides.Add(13, 2)
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Employees are required to report financial discrepancies immediately`rMaternity leave is provided for 90 days, while paternity leave is available for 30 ."

$s = $p.Slides.Add(14, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Key Points"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Standard work hours are from 9 AM to 5 PM, Monday to Friday`rOvertime is compensated at x1.5 regular pay for extra hours"

$s = $p.Slides.Add(15, 2)
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Remote work is subject to approval based on role requirements`rCybersecurity protocols must be followed when accessing company systems remotely ."

$s = $p.Slides.Add(16, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Key Points"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Ethics Hotline or reporting system for employees to report unethical behavior or concerns without fear of retaliation`rDisciplinary actions may include warnings, suspension, or termination"

$s = $p.Slides.Add(17, 2)
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Open-Door Policy: Employees may approach HR, direct managers, or senior leadership to discuss ethical concerns without formal procedures ."

$s = $p.Slides.Add(18, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Key Points"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Tier 1 – Entry-Level (Interns, Junior Developers, Support Staff) Employees in this tier perform basic tasks such as coding, testing, and assisting senior staff`rTier 2 – Mid-Level employees are responsible for designing, coding, leading small teams, executing projects, and managing financial operations"

$s = $p.Slides.Add(19, 2)
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Tier 3 – Senior-level (Senior Developers, Managers, HR Heads)"

$s = $p.Slides.Add(20, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Key Points"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Tier 4 – Executive-Level (Directors, C-Suite, Board Members) requires strategic thinking, advanced management, and strong stakeholder communication skills`rThey must possess visionary leadership, corporate governance expertise, and  visionary leadership"

$s = $p.Slides.Add(21, 2)
$s.Shapes.Item(2).TextFrame.TextRange.Text = "The BitBlock policy framework will be reviewed by company stakeholders and human  resources within 3 weeks ."

$s = $p.Slides.Add(22, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Key Points"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Any revision in the proposed changes will undergo the same procedure from the  beginning`rThe changes will be implemented in the same manner from the beginning of the year 2025"

$s = $p.Slides.Add(23, 2)
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Any changes will have to be made to comply with the requirements of the current system ."

$s = $p.Slides.Add(24, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Summary"
$s.Shapes.Item(2).TextFrame.TextRange.Text = " Any revision in the proposed changes will undergo the same procedure from the  beginning . The changes will be implemented in the same manner from the beginning of the year 2025 . Any changes will have to be made to comply with the requirements of the current system ."

